$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph that originally followed the
#    H1 title ("Play 8 Dragons Slot Game for Free - Review").
$d.Paragraphs.Item(2).Range.Delete() | Out-Null

# 2. Insert a new paragraph (bold "Play 8 Dragons Slot Game for Free - Review")
#    right before the final paragraph (the italic image-prompt paragraph).
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
$lastPara.Range.InsertParagraphBefore() | Out-Null

$n2 = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($n2 - 1)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play 8 Dragons Slot Game for Free - Review</w:t></w:r></w:p>'
$newPara.Range.InsertXML($newParaXml) | Out-Null

# 3. Replace the text of the final (italic) paragraph with the meta-description
#    copy (dropping the leading ": ").
$d.Content.Find.Execute("Create a feature image for `"8 Dragons`" that features a happy Maya warrior with glasses. The image should be in a cartoon style and should have a vibrant and eye-catching color scheme. The Maya warrior should be depicted holding a dragon in one hand and a pile of gold coins in the other, surrounded by Chinese-themed symbols such as lanterns and scrolls. In the background, you can add a colorful dragon or a temple to add to the overall theme of the game. The image should convey the excitement and adventure of playing `"8 Dragons`" and entice players to try their luck at this exciting slot game.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of 8 Dragons, an online slot game with impeccable graphics and a chance to choose free spins and multipliers. Play for free today.", 2) | Out-Null
